$d = $word.ActiveDocument
$table = $d.Tables.Item(1)

# Row 1: 99.92 -> 0M
$table.Cell(1,1).Range.Text = "0M"

# Row 2: 0.02 -> 0M
$table.Cell(2,1).Range.Text = "0M"

# Row 3: 29 -> 0M
$table.Cell(3,1).Range.Text = "0M"

# Row 4: 71 -> 191
$table.Cell(4,1).Range.Text = "191"

# Row 5: 0.00004 -> 0.00002
$table.Cell(5,1).Range.Text = "0.00002"

# Row 6: 0.00027 -> 0.00047
$table.Cell(6,1).Range.Text = "0.00047"

# Row 7: 0.00009 -> 0.00011
$table.Cell(7,1).Range.Text = "0.00011"

# Row 9: 0.00010 -> 0.00017
$table.Cell(9,1).Range.Text = "0.00017"

# Row 10: 0.00011 -> 0.00018
$table.Cell(10,1).Range.Text = "0.00018"

# Row 11: 0.00012 -> 0.00021
$table.Cell(11,1).Range.Text = "0.00021"

# Row 12: 0.00802 -> 0.02383
$table.Cell(12,1).Range.Text = "0.02383"

# Row 44: collapses a tab-delimited run into a single value "99.92"
$table.Cell(44,1).Range.Text = "99.92"

# Row 45: collapses a tab-delimited run into a single value "0.02"
$table.Cell(45,1).Range.Text = "0.02"

# Row 46: collapses a tab-delimited run into a single value "29"
$table.Cell(46,1).Range.Text = "29"
